# Scheduled runner update: refresh market-board derived values (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job Profits worksheets.
# Generated from the recorded price-refresh diff; applies plain value writes only
# (source data has no formulas), clearing cells that became blank and creating
# cells that became populated so the sheet matches the refreshed snapshot exactly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 861.5
$ws.Range("I4").Value = 861.5
$ws.Range("K4").Value = 861.5
$ws.Range("M4").Value = -747.5
# Row 6
$ws.Range("H6").Value = 743.2857
$ws.Range("J6").Value = 780.2727
$ws.Range("L6").Value = 2340.8181
$ws.Range("N6").Value = -2564.8181
# Row 15
$ws.Range("H15").Value = 781.8193
$ws.Range("I15").Value = 781.8193
$ws.Range("K15").Value = 2345.4579
$ws.Range("M15").Value = -2176.4579
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
# Row 39
$ws.Range("H39").Value = 595.3333
$ws.Range("I39").Value = 455.625
$ws.Range("J39").Value = 874.75
$ws.Range("K39").Value = 1366.875
$ws.Range("L39").Value = 2624.25
$ws.Range("M39").Value = -1070.875
$ws.Range("N39").Value = -3216.25
# Row 81
$ws.Range("H81").Value = 72198.75
$ws.Range("J81").Value = 72198.75
$ws.Range("L81").Value = 72198.75
$ws.Range("N81").Value = -74194.75
# Row 84
$ws.Range("H84").Value = 72198.75
$ws.Range("J84").Value = 72198.75
$ws.Range("L84").Value = 216596.25
$ws.Range("N84").Value = -226580.25
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4642.4443
$ws.Range("I61").Value = 2356.8
$ws.Range("J61").Value = 7499.5
$ws.Range("K61").Value = 2356.8
$ws.Range("L61").Value = 7499.5
$ws.Range("M61").Value = -2144.8
$ws.Range("N61").Value = -7923.5
# Row 102
$ws.Range("H102").Value = 4624.75
$ws.Range("I102").Value = 4624.75
$ws.Range("K102").Value = 4624.75
$ws.Range("M102").Value = -3002.75
# Row 132
$ws.Range("H132").Value = 3798.3635
$ws.Range("I132").Value = 4032.75
$ws.Range("J132").Value = 3664.4285
$ws.Range("K132").Value = 12098.25
$ws.Range("L132").Value = 10993.2855
$ws.Range("M132").Value = -9568.25
$ws.Range("N132").Value = -16053.2855
# Row 136
$ws.Range("H136").Value = 4642.4443
$ws.Range("I136").Value = 2356.8
$ws.Range("J136").Value = 7499.5
$ws.Range("K136").Value = 7070.400000000001
$ws.Range("L136").Value = 22498.5
$ws.Range("M136").Value = -4520.400000000001
$ws.Range("N136").Value = -27598.5
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1302.8334
$ws.Range("I86").Value = 1373.2307
$ws.Range("K86").Value = 1373.2307
$ws.Range("M86").Value = -250.2307000000001
# Row 89
$ws.Range("H89").Value = 1302.8334
$ws.Range("I89").Value = 1373.2307
$ws.Range("K89").Value = 6866.1535
$ws.Range("M89").Value = -1250.1535
# Row 99
$ws.Range("H99").Value = 2643.125
$ws.Range("I99").Value = 1488.2222
$ws.Range("J99").Value = 4128
$ws.Range("K99").Value = 1488.2222
$ws.Range("L99").Value = 4128
$ws.Range("M99").Value = 9.77780000000007
$ws.Range("N99").Value = -7124
# Row 134
$ws.Range("H134").Value = 3815.6875
$ws.Range("I134").Value = 3815.6875
$ws.Range("K134").Value = 11447.0625
$ws.Range("M134").Value = -8912.0625
$ws = $wb.Worksheets.Item("CRP")
# Row 51
$ws.Range("H51").Value = 31172.666
$ws.Range("J51").Value = 31172.666
$ws.Range("L51").Value = 31172.666
$ws.Range("N51").Value = -32644.666
# Row 58
$ws.Range("H58").Value = 2968.4285
$ws.Range("I58").Value = 2723.5
$ws.Range("J58").Value = 3295
$ws.Range("K58").Value = 2723.5
$ws.Range("L58").Value = 3295
$ws.Range("M58").Value = -2520.5
$ws.Range("N58").Value = -3701
# Row 61
$ws.Range("H61").Value = 31172.666
$ws.Range("J61").Value = 31172.666
$ws.Range("L61").Value = 31172.666
$ws.Range("N61").Value = -31868.666
# Row 74
$ws.Range("H74").Value = 66666.336
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19126
# Row 77
$ws.Range("H77").Value = 66666.336
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55632
# Row 132
$ws.Range("H132").Value = 88329.86
$ws.Range("I132").Value = 151627.5
$ws.Range("K132").Value = 454882.5
$ws.Range("M132").Value = -452352.5
# Row 134
$ws.Range("H134").Value = 2373.6365
$ws.Range("I134").Value = 2061
$ws.Range("K134").Value = 6183
$ws.Range("M134").Value = -3648
# Row 136
$ws.Range("H136").Value = 2968.4285
$ws.Range("I136").Value = 2723.5
$ws.Range("J136").Value = 3295
$ws.Range("K136").Value = 8170.5
$ws.Range("L136").Value = 9885
$ws.Range("M136").Value = -5620.5
$ws.Range("N136").Value = -14985
$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 1840.6072
$ws.Range("I132").Value = 1738
$ws.Range("J132").Value = 1917.5625
$ws.Range("K132").Value = 15642
$ws.Range("L132").Value = 17258.0625
$ws.Range("M132").Value = -13112
$ws.Range("N132").Value = -22318.0625
$ws = $wb.Worksheets.Item("GSM")
# Row 125
$ws.Range("H125").Value = 81857.836
$ws.Range("J125").Value = 81857.836
$ws.Range("L125").Value = 81857.836
$ws.Range("N125").Value = -86777.836
# Row 132
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 6000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -11060
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 5430.9473
$ws.Range("I46").Value = 1772.1428
$ws.Range("J46").Value = 7565.25
$ws.Range("K46").Value = 1772.1428
$ws.Range("L46").Value = 7565.25
$ws.Range("M46").Value = -1584.1428
$ws.Range("N46").Value = -7941.25
# Row 55
$ws.Range("H55").Value = 2715.3125
$ws.Range("I55").Value = 355.2
$ws.Range("K55").Value = 355.2
$ws.Range("M55").Value = -182.2
# Row 61
$ws.Range("H61").Value = 1246.25
$ws.Range("I61").Value = 1102.8572
$ws.Range("K61").Value = 1102.8572
$ws.Range("M61").Value = -900.8571999999999
# Row 81
$ws.Range("H81").Value = 34181
$ws.Range("J81").Value = 34181
$ws.Range("L81").Value = 34181
$ws.Range("N81").Value = -36177
# Row 82
$ws.Range("H82").Value = 11320.75
$ws.Range("I82").Value = 10720.643
$ws.Range("K82").Value = 10720.643
$ws.Range("M82").Value = -10359.643
# Row 84
$ws.Range("H84").Value = 34181
$ws.Range("J84").Value = 34181
$ws.Range("L84").Value = 102543
$ws.Range("N84").Value = -112527
# Row 85
$ws.Range("H85").Value = 11320.75
$ws.Range("I85").Value = 10720.643
$ws.Range("K85").Value = 10720.643
$ws.Range("M85").Value = -9472.643
# Row 100
$ws.Range("H100").Value = 2942.0527
$ws.Range("I100").Value = 2812.4375
$ws.Range("K100").Value = 2812.4375
$ws.Range("M100").Value = -2271.4375
# Row 113
$ws.Range("H113").Value = 1246.25
$ws.Range("I113").Value = 1102.8572
$ws.Range("K113").Value = 1102.8572
$ws.Range("M113").Value = 1067.1428
# Row 116
$ws.Range("H116").Value = 100000
$ws.Range("J116").Value = 100000
$ws.Range("L116").Value = 100000
$ws.Range("N116").Value = -109178
# Row 132
$ws.Range("H132").Value = 6798.231
$ws.Range("I132").Value = 6798.231
$ws.Range("K132").Value = 20394.693
$ws.Range("M132").Value = -17864.693
$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
# Row 96
$ws.Range("H96").Value = 2214.4
$ws.Range("I96").Value = 1849.75
$ws.Range("J96").Value = 2457.5
$ws.Range("K96").Value = 1849.75
$ws.Range("L96").Value = 2457.5
$ws.Range("M96").Value = -476.75
$ws.Range("N96").Value = -5203.5
# Row 100
$ws.Range("H100").Value = 812.25
$ws.Range("I100").Value = 812.25
$ws.Range("K100").Value = 1624.5
$ws.Range("M100").Value = -1083.5
# Row 132
$ws.Range("H132").Value = 4072.0967
$ws.Range("I132").Value = 4043.5
$ws.Range("K132").Value = 12130.5
$ws.Range("M132").Value = -9600.5
